$wb = $excel.ActiveWorkbook

function Add-Rows {
    param($ws, [int]$startRow, $rows)
    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $startRow + $i
        $row = $rows[$i]
        for ($c = 0; $c -lt $row.Count; $c++) {
            $ws.Cells.Item($r, $c + 1).Value = $row[$c]
        }
    }
}

# ALERTS: add rows starting at 23
$rows_ALERTS = @(
    ,@('''2026-01-30', '18:11:14', '18:00', 'Living Room', 'CRITICAL', 'FALL_DETECTED')
    ,@('''2026-01-30', '18:11:18', '18:00', 'Living Room', 'CRITICAL', 'FALL_DETECTED')
)
$ws = $wb.Worksheets.Item("ALERTS")
Add-Rows $ws 23 $rows_ALERTS

# PIR: add rows starting at 441
$rows_PIR = @(
    ,@('''2026-01-30', '18:07:59', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:08:02', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:08:02', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:08:07', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:08:12', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:08:17', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:08:23', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:08:27', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:08:32', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:08:37', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:08:42', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:11:19', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:11:22', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:11:23', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:11:28', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:11:33', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:11:39', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:11:43', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:11:48', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:11:53', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:11:59', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:12:03', '18:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('''2026-01-30', '18:12:08', '18:00', 'Bathroom', 'No Motion', 'Inactive')
)
$ws = $wb.Worksheets.Item("PIR")
Add-Rows $ws 441 $rows_PIR

# Humidity: add rows starting at 289
$rows_Humidity = @(
    ,@('''2026-01-30', '18:08:03', '18:00', 'Bathroom', '''86.8%', 'Active')
    ,@('''2026-01-30', '18:08:13', '18:00', 'Bathroom', '''86.8%', 'Active')
    ,@('''2026-01-30', '18:08:24', '18:00', 'Bathroom', '''86.7%', 'Active')
    ,@('''2026-01-30', '18:08:33', '18:00', 'Bathroom', '''86.7%', 'Active')
    ,@('''2026-01-30', '18:08:43', '18:00', 'Bathroom', '''86.7%', 'Active')
    ,@('''2026-01-30', '18:11:20', '18:00', 'Bathroom', '''85.8%', 'Active')
    ,@('''2026-01-30', '18:11:24', '18:00', 'Bathroom', '''85.9%', 'Active')
    ,@('''2026-01-30', '18:11:29', '18:00', 'Bathroom', '''86.8%', 'Active')
    ,@('''2026-01-30', '18:11:34', '18:00', 'Bathroom', '''85.8%', 'Active')
    ,@('''2026-01-30', '18:11:40', '18:00', 'Bathroom', '''86.8%', 'Active')
    ,@('''2026-01-30', '18:11:44', '18:00', 'Bathroom', '''86.9%', 'Active')
    ,@('''2026-01-30', '18:11:49', '18:00', 'Bathroom', '''86.9%', 'Active')
    ,@('''2026-01-30', '18:11:54', '18:00', 'Bathroom', '''86.1%', 'Active')
    ,@('''2026-01-30', '18:12:00', '18:00', 'Bathroom', '''87.3%', 'Active')
    ,@('''2026-01-30', '18:12:04', '18:00', 'Bathroom', '''86.2%', 'Active')
    ,@('''2026-01-30', '18:12:10', '18:00', 'Bathroom', '''86.9%', 'Active')
)
$ws = $wb.Worksheets.Item("Humidity")
Add-Rows $ws 289 $rows_Humidity

# Proximity: add rows starting at 97
$rows_Proximity = @(
    ,@('''2026-01-30', '18:07:58', '18:00', 'Living Room Main Door', 'ENTER', 'User ENTERED Living Room Main Door')
    ,@('''2026-01-30', '18:08:01', '18:00', 'Living Room Main Door', 'EXIT', 'User EXITED Living Room Main Door')
    ,@('''2026-01-30', '18:08:14', '18:00', 'Living Room Main Door', 'ENTER', 'User ENTERED Living Room Main Door')
    ,@('''2026-01-30', '18:08:22', '18:00', 'Living Room Main Door', 'EXIT', 'User EXITED Living Room Main Door')
    ,@('''2026-01-30', '18:08:29', '18:00', 'Living Room Main Door', 'ENTER', 'User ENTERED Living Room Main Door')
    ,@('''2026-01-30', '18:11:21', '18:00', 'Living Room Main Door', 'ENTER', 'User ENTERED Living Room Main Door')
    ,@('''2026-01-30', '18:11:38', '18:00', 'Living Room Main Door', 'EXIT', 'User EXITED Living Room Main Door')
    ,@('''2026-01-30', '18:11:58', '18:00', 'Living Room Main Door', 'ENTER', 'User ENTERED Living Room Main Door')
)
$ws = $wb.Worksheets.Item("Proximity")
Add-Rows $ws 97 $rows_Proximity

# Camera: add rows starting at 24
$rows_Camera = @(
    ,@('''2026-01-30', '18:07:57', '18:00', 'Living Room Main Door', 'Image Captured (ENTER)', 'Active')
    ,@('''2026-01-30', '18:08:00', '18:00', 'Living Room Main Door', 'Image Captured (EXIT)', 'Active')
    ,@('''2026-01-30', '18:08:14', '18:00', 'Living Room Main Door', 'Image Captured (ENTER)', 'Active')
    ,@('''2026-01-30', '18:08:21', '18:00', 'Living Room Main Door', 'Image Captured (EXIT)', 'Active')
    ,@('''2026-01-30', '18:08:28', '18:00', 'Living Room Main Door', 'Image Captured (ENTER)', 'Active')
)
$ws = $wb.Worksheets.Item("Camera")
Add-Rows $ws 24 $rows_Camera

